# MAJ fonction informations compte / Planning / UserStories
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update progress values (column B)
$ws.Range("B5").Value = 0.8
$ws.Range("B6").Value = 0.9
$ws.Range("B8").Value = 0.8
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 0.8
$ws.Range("B11").Value = 0.9
$ws.Range("B12").Value = 0.9
$ws.Range("B14").Value = 0.7
$ws.Range("B15").Value = 0.3
$ws.Range("B24").Value = 0.7

# Mark "X" in column H for rows 10, 11, 12, 14
$ws.Range("H10").Value = "X"
$ws.Range("H11").Value = "X"
$ws.Range("H12").Value = "X"
$ws.Range("H14").Value = "X"

# Update active selection to B25
$ws.Range("B25").Select()
